$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2242.182
$ws.Range("I15").Value = 2242.182
$ws.Range("K15").Value = 6726.545999999999
$ws.Range("M15").Value = -6557.545999999999

$ws.Range("H43").Value = 1848.9286
$ws.Range("I43").Value = 700
$ws.Range("J43").Value = 2162.2727
$ws.Range("K43").Value = 700
$ws.Range("L43").Value = 2162.2727
$ws.Range("M43").Value = -631
$ws.Range("N43").Value = -2300.2727

$ws.Range("H116").Value = 3362.7026
$ws.Range("I116").Value = 2501.2083
$ws.Range("K116").Value = 2501.2083
$ws.Range("M116").Value = 940.7917000000002

$ws.Range("H132").Value = 8701752
$ws.Range("I132").Value = 9529538
$ws.Range("K132").Value = 28588614
$ws.Range("M132").Value = -28586084

$ws.Range("H135").Value = 789.05
$ws.Range("I135").Value = 820.0526
$ws.Range("K135").Value = 7380.4734
$ws.Range("M135").Value = -4845.4734

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5401.6606
$ws.Range("I32").Value = 3791.625
$ws.Range("J32").Value = 15061.875
$ws.Range("K32").Value = 3791.625
$ws.Range("L32").Value = 15061.875
$ws.Range("M32").Value = -3504.625
$ws.Range("N32").Value = -15635.875

$ws.Range("H45").Value = 1844.4073
$ws.Range("I45").Value = 1104.125
$ws.Range("K45").Value = 1104.125
$ws.Range("M45").Value = -727.125

$ws.Range("H110").Value = 1112.8823
$ws.Range("I110").Value = 530.25806
$ws.Range("K110").Value = 530.25806
$ws.Range("M110").Value = 1514.74194

$ws.Range("H122").Value = 3142.7896
$ws.Range("I122").Value = 2181.3
$ws.Range("J122").Value = 4211.1113
$ws.Range("K122").Value = 6543.900000000001
$ws.Range("L122").Value = 12633.3339
$ws.Range("M122").Value = -4093.900000000001
$ws.Range("N122").Value = -17533.3339

$ws.Range("H132").Value = 21280008
$ws.Range("I132").Value = 31253068
$ws.Range("J132").Value = 4148.4
$ws.Range("K132").Value = 93759204
$ws.Range("L132").Value = 12445.2
$ws.Range("M132").Value = -93756674
$ws.Range("N132").Value = -17505.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1542.25
$ws.Range("I105").Value = 1207.3334
$ws.Range("J105").Value = 2100.4443
$ws.Range("K105").Value = 1207.3334
$ws.Range("L105").Value = 2100.4443
$ws.Range("M105").Value = 539.6666
$ws.Range("N105").Value = -5594.4443

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3125.1667
$ws.Range("I16").Value = 3097.8
$ws.Range("J16").Value = 3144.7144
$ws.Range("K16").Value = 3097.8
$ws.Range("L16").Value = 3144.7144
$ws.Range("M16").Value = -2810.8
$ws.Range("N16").Value = -3718.7144

$ws.Range("H31").Value = 2555.3208
$ws.Range("I31").Value = 1818.4889
$ws.Range("J31").Value = 6700
$ws.Range("K31").Value = 1818.4889
$ws.Range("L31").Value = 6700
$ws.Range("M31").Value = -1523.4889
$ws.Range("N31").Value = -7290

$ws.Range("H34").Value = 2555.3208
$ws.Range("I34").Value = 1818.4889
$ws.Range("J34").Value = 6700
$ws.Range("K34").Value = 1818.4889
$ws.Range("L34").Value = 6700
$ws.Range("M34").Value = -1616.4889
$ws.Range("N34").Value = -7104

$ws.Range("H113").Value = 3125.1667
$ws.Range("I113").Value = 3097.8
$ws.Range("J113").Value = 3144.7144
$ws.Range("K113").Value = 3097.8
$ws.Range("L113").Value = 3144.7144
$ws.Range("M113").Value = -927.8000000000002
$ws.Range("N113").Value = -7484.7144

$ws.Range("H132").Value = 3522.682
$ws.Range("I132").Value = 2353.4167
$ws.Range("K132").Value = 7060.250100000001
$ws.Range("M132").Value = -4530.250100000001

$ws.Range("H134").Value = 1379.4736
$ws.Range("I134").Value = 706.4483
$ws.Range("J134").Value = 3548.111
$ws.Range("K134").Value = 2119.3449
$ws.Range("L134").Value = 10644.333
$ws.Range("M134").Value = 415.6550999999999
$ws.Range("N134").Value = -15714.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1637.9
$ws.Range("I4").Value = 87.5
$ws.Range("J4").Value = 2025.5
$ws.Range("K4").Value = 262.5
$ws.Range("L4").Value = 6076.5
$ws.Range("M4").Value = -150.5
$ws.Range("N4").Value = -6300.5

$ws.Range("H122").Value = 1948.1
$ws.Range("I122").Value = 931.3333
$ws.Range("J122").Value = 2383.8572
$ws.Range("K122").Value = 8381.9997
$ws.Range("L122").Value = 21454.7148
$ws.Range("M122").Value = -5931.9997
$ws.Range("N122").Value = -26354.7148

$ws.Range("H131").Value = 1060.0571
$ws.Range("I131").Value = 546.6667
$ws.Range("J131").Value = 1166.2759
$ws.Range("K131").Value = 1640.0001
$ws.Range("L131").Value = 3498.8277
$ws.Range("M131").Value = 3399.9999
$ws.Range("N131").Value = -13578.8277

$ws.Range("H136").Value = 2087.7646
$ws.Range("J136").Value = 3806.6
$ws.Range("L136").Value = 11419.8
$ws.Range("N136").Value = -21619.8

$ws.Range("H138").Value = 1408.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3375
$ws.Range("I80").Value = 3375
$ws.Range("K80").Value = 3375
$ws.Range("M80").Value = -2377

$ws.Range("H83").Value = 3375
$ws.Range("I83").Value = 3375
$ws.Range("K83").Value = 16875
$ws.Range("M83").Value = -11883

$ws.Range("H132").Value = 3107.75
$ws.Range("I132").Value = 2565.762
$ws.Range("K132").Value = 7697.286
$ws.Range("M132").Value = -5167.286

$ws.Range("H137").Value = 29618.75
$ws.Range("J137").Value = 29618.75
$ws.Range("L137").Value = 29618.75
$ws.Range("N137").Value = -39818.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 142859470
$ws.Range("I22").Value = 333333800
$ws.Range("J22").Value = 3735.5
$ws.Range("K22").Value = 333333800
$ws.Range("L22").Value = 3735.5
$ws.Range("M22").Value = -333333505
$ws.Range("N22").Value = -4325.5

$ws.Range("H27").Value = 142859470
$ws.Range("I27").Value = 333333800
$ws.Range("J27").Value = 3735.5
$ws.Range("K27").Value = 333333800
$ws.Range("L27").Value = 3735.5
$ws.Range("M27").Value = -333333693
$ws.Range("N27").Value = -3949.5

$ws.Range("H40").Value = 2413.9
$ws.Range("I40").Value = 1200
$ws.Range("J40").Value = 2934.1428
$ws.Range("K40").Value = 1200
$ws.Range("L40").Value = 2934.1428
$ws.Range("M40").Value = -1064
$ws.Range("N40").Value = -3206.1428

$ws.Range("H140").Value = 29489.857
$ws.Range("J140").Value = 29489.857
$ws.Range("L140").Value = 29489.857
$ws.Range("N140").Value = -39849.857

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 346544.06
$ws.Range("I122").Value = 456055.22
$ws.Range("K122").Value = 1368165.66
$ws.Range("M122").Value = -1365715.66
